# Apply updated crypto price/volume data as produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.006.26'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '3.432.54'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '571.61'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '159.13'
$ws.Range('E6').Value = '  +2.07%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.434.00'
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').Value = '0.571'
$ws.Range('E9').Value = '  -8.78%  '
$ws.Range('D10').Value = '7.27'
$ws.Range('E10').Value = '  +1.96%  '
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('D12').Value = '0.424'
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('D13').Value = '4.013.02'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '0.135'
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '27.18'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('D16').Value = '''0.0000173'
$ws.Range('E16').Value = '  -7.20%  '
$ws.Range('D17').Value = '64.037.07'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '3.411.41'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').Value = '''6.10'
$ws.Range('E19').Value = '  -3.12%  '
$ws.Range('D20').Value = '13.66'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '384.34'
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('D22').Value = '7.87'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').Value = '''71.30'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '0.516'
$ws.Range('E25').Value = '  -5.07%  '
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('D27').Value = '9.69'
$ws.Range('E27').Value = '  -5.79%  '
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('D32').Value = '1.99'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '22.99'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '6.98'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').Value = '1.52'
$ws.Range('E36').Value = '  -5.35%  '
$ws.Range('D37').Value = '''160.90'
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('D38').Value = '0.849'
$ws.Range('E38').Value = '  +11.00%  '
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').Value = '  -2.56%  '
$ws.Range('D40').Value = '2.827.63'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').Value = '26.07'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').Value = '0.0725'
$ws.Range('E42').Value = '  -4.59%  '
$ws.Range('D43').Value = '43.04'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('D45').Value = '6.42'
$ws.Range('E45').Value = '  -7.76%  '
$ws.Range('D46').Value = '4.37'
$ws.Range('E46').Value = '  -5.07%  '
$ws.Range('D47').Value = '0.0305'
$ws.Range('E47').Value = '  -2.94%  '
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +11.27%  '
$ws.Range('D49').Value = '333.58'
$ws.Range('E49').Value = '  +4.30%  '
$ws.Range('D50').Value = '1.05'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('E51').Value = '  -5.23%  '
